# Update the document date and the multiplication problem/answer pairs
# in the table, per the commit "Update master to output generated at
# c8c62b6".

$d = $word.ActiveDocument

# Map of old text -> new text. Each old string is unique within the
# document, so a simple Find/Replace (wdReplaceAll) per pair is safe.
$replacements = [ordered]@{
    "2025-06-17 Tuesday" = "2025-06-18 Wednesday"
    "372×2=744"           = "601×5=3005"
    "735×6=4410"          = "756×8=6048"
    "121×6=726"           = "693×4=2772"
    "354×2=708"           = "884×8=7072"
    "854×7=5978"          = "958×3=2874"
    "335×8=2680"          = "158×7=1106"
    "817×8=6536"          = "830×3=2490"
    "313×2=626"           = "686×9=6174"
    "627×2=1254"          = "602×3=1806"
    "804×6=4824"          = "738×9=6642"
    "296×9=2664"          = "187×9=1683"
    "814×9=7326"          = "378×9=3402"
    "175×9=1575"          = "265×3=795"
    "285×4=1140"          = "685×7=4795"
    "226×2=452"           = "388×7=2716"
    "944×2=1888"          = "810×9=7290"
    "933×7=6531"          = "885×5=4425"
    "128×8=1024"          = "190×6=1140"
    "922×9=8298"          = "149×6=894"
    "697×2=1394"          = "124×8=992"
    "516×5=2580"          = "150×3=450"
    "442×3=1326"          = "461×2=922"
    "319×4=1276"          = "833×4=3332"
    "610×5=3050"          = "389×7=2723"
    "865×3=2595"          = "610×7=4270"
}

foreach ($old in $replacements.Keys) {
    $new = $replacements[$old]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new, 2)
}

$d.Save()
